# Development Log - add a new entry at the top of the log table (row 4),
# pushing the existing entries down by one row.
#
# Commit message:
#   Fix logic using while loop to continuously prompt user, if they do not
#   enter the correct format. - i.e., word followed by comma, and letter
#   (or _) and up to two numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev Log")

# --- Insert a new blank row above the current row 4 (the table's first
#     data row), shifting the existing entries down --------------------
$ws.Rows.Item(4).Insert()

# Match the formatting of the row immediately below (which is the entry
# that used to be row 4), so the new row renders identically to the rest
# of the log table.
$ws.Range("A5:G5").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = 298

# --- Populate the new log entry ----------------------------------------
$ws.Range("B4").Value = 45364
$ws.Range("C4").Value = 0.22708333333333333
$ws.Range("D4").Value = "enterWordAndTile()"
$ws.Range("F4").Value = "Fixing logic using while loop to continuously prompt user, if they do not enter the correct format. - ine.m word, letter (or _) and up to two numbers."
$ws.Range("E4").Value = "'trying out different loops and string formatting techniques"
$ws.Range("G4").Value = 0.56000000000000005

# --- Misc view tweaks recorded in the same commit ----------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("A5").Select()

Write-Host "Added new Dev Log entry (enterWordAndTile) at row 4"
